$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing shared string (si index 60 / cell E15 'restructuring' text): trim trailing status annotations ---
$ws.Range("E15").Value = "Each script (scraper, cleaner, etc.) was previously built in an ad-hoc, local way with minimal modularity and no clear interface for
orchestration. In order to integrate them with Prefect and GitHub Actions, I had to:
Modularize the scripts to expose main functions
     Improve logging and error handling
     Add configuration via .env or parameters for flexibility
     Create reusable Prefect flows and GitHub workflows that can be scheduled and monitored
     Reproduce the same behavior and logs from the CLI when triggered through GitHub
This restructuring also had to preserve all existing features like:
     Extracting from multiple subreddits
     Next page has no posts from yesterday
     Next page has mixed days, only collect yesterday   
     Less than 500 posts available  "

# --- Row 15: update Status (F15) text + enable wrap ---
$ws.Range("F15").WrapText = $true
$ws.Range("F15").Value = "Partially Completed: The idempotency results or the main features and their idempotency errors needs extra checking
Condition	                                                            Status	               Notes
Extract from multiple subreddits	                               ✅ Fully Met	No issues.
Stop when next page has no posts from yesterday	🟡 Partial	You filter, but don’t stop early. Add a check for all irrelevant posts.
Only collect posts from yesterday if mixed days	 ✅ Fully Met	Correct filtering.
Accept fewer than 500 posts, no older-day spillover	🟡 Partial	Filtering works, but no early exit when all new posts are too old."

# --- Row 16 ---
$ws.Range("A16").Value = "23/7/2025(Onsite)"
$ws.Range("B16").Value = "Car Tracking Project"
$ws.Range("C16").WrapText = $true
$ws.Range("C16").Value = "Add debugging counters so that you can track the progress and the pressure that you apply on reddit's api"
$ws.Range("D16").WrapText = $true
$ws.Range("F16").Value = "DONE"
$ws.Rows.Item(16).RowHeight = 28.8

# --- Row 17 ---
$ws.Range("A17").Value = "23/7/2025(Onsite)"
$ws.Range("B17").Value = "Car Tracking Project"
$ws.Range("C17").WrapText = $true
$ws.Range("C17").Value = "Test if the counters are working well when you change the t param from params dictionary from
 week to day"
$ws.Range("F17").WrapText = $true
$ws.Range("F17").Value = "DONE:
day param output with the selected columns:
=== Debugging Counters ===
total_posts_fetched: 11792
posts_filtered_time: 10855
posts_filtered_comments: 41
comments_skipped: 715
valid_posts_stored: 665
Completed in 15.40 minutes
week param output with the selected columns:
=== Debugging Counters ===
total_posts_fetched: 11875
posts_filtered_time: 10938
posts_filtered_comments: 41
comments_skipped: 715
valid_posts_stored: 666
Completed in 15.15 minutes
Not that big of a difference, need to understand the underlying logic of the built in functions more. Currently at this stage I care about getting the MVP. Next step: automate the extraction process using prefect and github actions, then find a way to clean the unstructured textual columns in the extracted csv files."
$ws.Rows.Item(17).RowHeight = 331.2

# --- Row 18 ---
$ws.Range("A18").Value = "24/7/2025(Vacation)"
$ws.Range("B18").Value = "Car Tracking Project"
$ws.Range("C18").Value = "Modularize the Extractor code so that you can scale the project later on"

# --- Row 19 ---
$ws.Range("A19").Value = "25/7/2025(Vacation)"
$ws.Range("B19").Value = "Car Tracking Project"
$ws.Range("C19").WrapText = $true
$ws.Range("C19").Value = "Make the script that will feed the data into The offline AI model and get its response and store it in 
CSV file"
$ws.Rows.Item(19).RowHeight = 28.8

# --- Column F width ---
$ws.Columns.Item(6).ColumnWidth = 85.67

# --- Extend used-range dimension artifact to CR19 (matches source file's dimension metadata) ---
$ws.Cells.Item(19, 96).Style = "Normal"

# --- Selection / view state ---
$ws.Range("A20").Select()

Write-Host "edit complete"
